$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new note row describing the CSV-like header convention (added to the
# shared string table first so it lands before the renamed headers below)
$ws.Range("A23").Value = "nome | tipo | unidade_dose | dose_minima | dose_maxima | ..."

# Update header row (C1:H1): replace spaces with underscores in the labels
$ws.Range("C1").Value = "Unidade_dose"
$ws.Range("D1").Value = "Dose_mínima"
$ws.Range("E1").Value = "Dose_máxima"
$ws.Range("F1").Value = "Concentração_máxima"
$ws.Range("G1").Value = "Diluição_sugerida"
$ws.Range("H1").Value = "Forma_de_administração"

# Update the selection to mirror the saved view
$ws.Range("J16").Select()
